$p = $ppt.ActivePresentation

# Slide 6: "Data Cleansing and Exploration" - body placeholder (shape 2)
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1,1).Runs(1)
$tr6.Text = "Decided to drop 800 rows where the variable delinquency and derogatory reposrts were missing."
$s6.Shapes.Item(2).TextFrame.AutoSize = 2

# Slide 7: "Data Cleansing and Exploration cont." - body placeholder (shape 2)
$s7 = $p.Slides.Item(7)
$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange.Paragraphs(2,1).Runs(1)
$tr7.Text = "Totaled number of rows with imputed data and categorical replaced was 5175."
